$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = [double]"2.0979733705356881E-2"
$ws.Range("B1").Value = [double]"1.9334749735755519E-3"
$ws.Range("C1").Value = [double]"0.27116203194303995"
$ws.Range("D1").Value = [double]"3.4231495030716527E-2"
$ws.Range("E1").Value = [double]"4.2292285895141498E-21"
$ws.Range("F1").Value = [double]"7.2101025099163912E-2"
$ws.Range("A2").Value = [double]"2.0968907239077843E-2"
$ws.Range("B2").Value = [double]"7.5525210371449551E-3"
$ws.Range("C2").Value = [double]"0.34441167593843824"
$ws.Range("D2").Value = [double]"4.0616936191589348E-2"
$ws.Range("E2").Value = [double]"9.5000000014450517E-21"
$ws.Range("F2").Value = [double]"0.70430158332812509"
$ws.Range("A3").Value = [double]"2.0858150336375655E-2"
$ws.Range("B3").Value = [double]"3.6907722091441401E-3"
$ws.Range("C3").Value = [double]"0.31535944071532468"
$ws.Range("D3").Value = [double]"3.5937457462386115E-2"
$ws.Range("E3").Value = [double]"9.5000000014450517E-21"
$ws.Range("F3").Value = [double]"0.22095445547354753"
$ws.Range("A4").Value = [double]"2.0441945242133948E-2"
$ws.Range("B4").Value = [double]"0.15709765275003271"
$ws.Range("C4").Value = [double]"0.23466073541094523"
$ws.Range("D4").Value = [double]"9.9999999999977801E-2"
$ws.Range("E4").Value = [double]"9.5001011555012064E-21"
$ws.Range("F4").Value = [double]"9.9999999999999787"
$ws.Range("A5").Value = [double]"2.0622658413166028E-2"
$ws.Range("B5").Value = [double]"9.1301546643056314E-2"
$ws.Range("C5").Value = [double]"0.27422375791725045"
$ws.Range("D5").Value = [double]"9.9999999999977801E-2"
$ws.Range("E5").Value = [double]"9.5000000014450517E-21"
$ws.Range("F5").Value = [double]"9.9999999999999787"
$ws.Range("A6").Value = [double]"2.0642900997480545E-2"
$ws.Range("B6").Value = [double]"0.17597151531869565"
$ws.Range("C6").Value = [double]"0.18059606056527022"
$ws.Range("D6").Value = [double]"9.9998912654022984E-2"
$ws.Range("E6").Value = [double]"9.5000000014450517E-21"
$ws.Range("F6").Value = [double]"7.4193704393690227"
$ws.Range("A7").Value = [double]"2.0840639164015299E-2"
$ws.Range("B7").Value = [double]"1.0333579136578971E-2"
$ws.Range("C7").Value = [double]"0.24848209941194072"
$ws.Range("D7").Value = [double]"3.8244186571606623E-2"
$ws.Range("E7").Value = [double]"5.5333888176577347E-21"
$ws.Range("F7").Value = [double]"8.7445039482058479E-2"
$ws.Range("A8").Value = [double]"2.0828526396279294E-2"
$ws.Range("B8").Value = [double]"1.8755404835566758E-2"
$ws.Range("C8").Value = [double]"0.29519935270273545"
$ws.Range("D8").Value = [double]"4.1436854877885961E-2"
$ws.Range("E8").Value = [double]"6.7396318658010479E-21"
$ws.Range("F8").Value = [double]"0.27852322066562413"
$ws.Range("A9").Value = [double]"2.0652503602074594E-2"
$ws.Range("B9").Value = [double]"7.5098722005191104E-2"
$ws.Range("C9").Value = [double]"0.28343621395806362"
$ws.Range("D9").Value = [double]"9.9999999999977801E-2"
$ws.Range("E9").Value = [double]"9.5000000014450517E-21"
$ws.Range("F9").Value = [double]"9.9999999999999787"
$ws.Range("A10").Value = [double]"2.0649151377202017E-2"
$ws.Range("B10").Value = [double]"7.7194533978412735E-2"
$ws.Range("C10").Value = [double]"0.14830011480480246"
$ws.Range("D10").Value = [double]"9.9999999999955777E-2"
$ws.Range("E10").Value = [double]"9.5000002601914286E-21"
$ws.Range("F10").Value = [double]"9.9999999999999574"
$ws.Range("A11").Value = [double]"2.0864682978647586E-2"
$ws.Range("B11").Value = [double]"1.0151935112894476E-2"
$ws.Range("C11").Value = [double]"0.35336944410388599"
$ws.Range("D11").Value = [double]"3.6977232911345373E-2"
$ws.Range("E11").Value = [double]"9.5000000014450517E-21"
$ws.Range("F11").Value = [double]"0.10329491146538314"
$ws.Range("A12").Value = [double]"1.9652762510853045E-2"
$ws.Range("B12").Value = [double]"0.37536021524588786"
$ws.Range("C12").Value = [double]"1.2384371789244055E-2"
$ws.Range("D12").Value = [double]"9.9999999999977801E-2"
$ws.Range("E12").Value = [double]"9.5000000014450517E-21"
$ws.Range("F12").Value = [double]"9.9999999999999787"
$ws.Range("A13").Value = [double]"2.0867267378827999E-2"
$ws.Range("B13").Value = [double]"1.0484980674986922E-2"
$ws.Range("C13").Value = [double]"0.34118378501272789"
$ws.Range("D13").Value = [double]"4.0036949983842947E-2"
$ws.Range("E13").Value = [double]"9.5019044347325565E-21"
$ws.Range("F13").Value = [double]"0.26698278381222834"
$ws.Range("A14").Value = [double]"2.0871816066514948E-2"
$ws.Range("B14").Value = [double]"4.8797271838518178E-3"
$ws.Range("C14").Value = [double]"0.28998110587356096"
$ws.Range("D14").Value = [double]"3.4755031777655321E-2"
$ws.Range("E14").Value = [double]"5.8857080544565737E-21"
$ws.Range("F14").Value = [double]"4.9925085065405161E-2"
$ws.Range("A15").Value = [double]"2.0791485710513324E-2"
$ws.Range("B15").Value = [double]"5.190239810621624E-2"
$ws.Range("C15").Value = [double]"0.37937061938419581"
$ws.Range("D15").Value = [double]"5.2980588291542487E-2"
$ws.Range("E15").Value = [double]"9.5000021566130448E-21"
$ws.Range("F15").Value = [double]"1.5654188917996501"
$ws.Range("A16").Value = [double]"2.0837783700202419E-2"
$ws.Range("B16").Value = [double]"2.1291611551281793E-2"
$ws.Range("C16").Value = [double]"0.34465048737378196"
$ws.Range("D16").Value = [double]"5.7655951953579851E-2"
$ws.Range("E16").Value = [double]"9.5000000014450517E-21"
$ws.Range("F16").Value = [double]"2.7598506753065086"

$wb.Save()
